$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values as literal text in the source data (e.g.
# "35.617.67", "1.00", "0.0608"). Force a Text number format on every Price cell
# we touch before writing so Excel keeps the exact digits/trailing zeros instead of
# silently coercing the string to a binary number (which would drop trailing zeros
# like "1.00" -> "1", or introduce floating-point noise like "0.694" -> "0.69399999999999995").
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D13", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D27", "D28", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D41", "D42", "D44", "D45", "D46", "D47", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "35.617.67"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.894.90"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "0.694"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "246.39"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "43.11"
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("D9").Value = "56.97"
$ws.Range("E9").Value = "  +9.30%  "
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "14.60"
$ws.Range("E13").Value = "  +10.38%  "
$ws.Range("E14").Value = "  +8.49%  "
$ws.Range("D15").Value = "2.170.93"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "5.05"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "1.905.69"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "35.580.54"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "73.61"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "0.0₃0831"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "246.26"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").Value = "13.01"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  +3.89%  "
$ws.Range("E24").Value = "  +4.76%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("D27").Value = "166.72"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "8.71"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "4.41"
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("D32").Value = "0.0608"
$ws.Range("E32").Value = "  +4.18%  "
$ws.Range("D33").Value = "4.27"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("B34").Value = "BinanceUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "1.84"
$ws.Range("E35").Value = "  +12.82%  "
$ws.Range("E36").Value = "  -17.43%  "
$ws.Range("D37").Value = "0.855"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.0740"
$ws.Range("E38").Value = "  +8.04%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "1.95"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("E40").Value = "  +6.77%  "
$ws.Range("D41").Value = "99.30"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "17.05"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").Value = "14.24"
$ws.Range("E44").Value = "  +17.36%  "
$ws.Range("D45").Value = "1.316.49"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "2.36"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").Value = "0.0809"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "6.41"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("E51").Value = "  -2.30%  "
